$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("RF")
$ws.Range("C2").Value = 0.7025693054032662
$ws.Range("D2").Value = 0.05543606843201483
$ws.Range("C3").Value = 0.8703512512624846
$ws.Range("D3").Value = 0.02373528245755575
$ws.Range("C4").Value = 0.9232680360089499
$ws.Range("D4").Value = 0.05040732636584883
$ws.Range("C5").Value = 0.5233763687293099
$ws.Range("D5").Value = 0.09329568105429881
$ws.Range("C6").Value = 0.9854839999999999
$ws.Range("D6").Value = 0.009426522721144038
$ws.Range("C7").Value = 0.6636276502512104
$ws.Range("D7").Value = 0.07867142953880846
$ws.Range("C8").Value = 0.8558668907595055
$ws.Range("D8").Value = 0.02947098778840414
$ws.Range("C9").Value = 0.7916072963737392
$ws.Range("D9").Value = 0.04601568035330826
$ws.Range("C10").Value = 0.754429368975313
$ws.Range("D10").Value = 0.04683649585922914
$ws.Range("C11").Value = 0.6306377336986995
$ws.Range("D11").Value = 0.07587842462280164
$ws.Range("C12").Value = 0.8622700000000001
$ws.Range("D12").Value = 0.02323600522412385
$ws.Range("C13").Value = 0.754429368975313
$ws.Range("D13").Value = 0.04683649585922915

$ws = $wb.Worksheets.Item("LGBM")
$ws.Range("C2").Value = 0.6902480971282161
$ws.Range("D2").Value = 0.05065605338441074
$ws.Range("C3").Value = 0.8803669621815732
$ws.Range("D3").Value = 0.0254031261764852
$ws.Range("C4").Value = 0.917351603102431
$ws.Range("D4").Value = 0.05155810095693884
$ws.Range("C5").Value = 0.5718070834394363
$ws.Range("D5").Value = 0.09091321025776683
$ws.Range("C6").Value = 0.9828880000000001
$ws.Range("D6").Value = 0.01060320783190306
$ws.Range("C7").Value = 0.7008704328115608
$ws.Range("D7").Value = 0.0757800303999466
$ws.Range("C8").Value = 0.8692345226512308
$ws.Range("D8").Value = 0.03016463394799619
$ws.Range("C9").Value = 0.8130062987750296
$ws.Range("D9").Value = 0.04528056898934198
$ws.Range("C10").Value = 0.7773466853498074
$ws.Range("D10").Value = 0.04666320125982015
$ws.Range("C11").Value = 0.6613456539497846
$ws.Range("D11").Value = 0.07750905278855336
$ws.Range("C12").Value = 0.874104
$ws.Range("D12").Value = 0.02408240072106176
$ws.Range("C13").Value = 0.7773466853498074
$ws.Range("D13").Value = 0.04666320125982015

$ws = $wb.Worksheets.Item("XGB")
$ws.Range("C2").Value = 0.7214059897225265
$ws.Range("D2").Value = 0.0532392516660679
$ws.Range("C3").Value = 0.8876848838514196
$ws.Range("D3").Value = 0.02593447855321011
$ws.Range("C4").Value = 0.9088277478061438
$ws.Range("D4").Value = 0.05669405501788487
$ws.Range("C5").Value = 0.6125289024700789
$ws.Range("D5").Value = 0.09511275833412092
$ws.Range("C6").Value = 0.979108
$ws.Range("D6").Value = 0.01379959833783992
$ws.Range("C7").Value = 0.7275418970355019
$ws.Range("D7").Value = 0.074489981792056
$ws.Range("C8").Value = 0.8788986172812676
$ws.Range("D8").Value = 0.03011668094895138
$ws.Range("C9").Value = 0.8283478791104728
$ws.Range("D9").Value = 0.04482266967727443
$ws.Range("C10").Value = 0.7958174164115089
$ws.Range("D10").Value = 0.04805861153055179
$ws.Range("C11").Value = 0.6836952594375814
$ws.Range("D11").Value = 0.07765960222975032
$ws.Range("C12").Value = 0.8844540000000001
$ws.Range("D12").Value = 0.02548602025761238
$ws.Range("C13").Value = 0.7958174164115089
$ws.Range("D13").Value = 0.04805861153055178

$ws = $wb.Worksheets.Item("KNN")
$ws.Range("C2").Value = 0.7041235999565607
$ws.Range("D2").Value = 0.06649094686717673
$ws.Range("C3").Value = 0.904703175850073
$ws.Range("D3").Value = 0.02045050563116494
$ws.Range("C4").Value = 0.9102850586764899
$ws.Range("D4").Value = 0.05464423477679217
$ws.Range("C5").Value = 0.6872465516509634
$ws.Range("D5").Value = 0.06988049368272724
$ws.Range("C6").Value = 0.9769239999999999
$ws.Range("D6").Value = 0.01437884727575864
$ws.Range("C7").Value = 0.7809453002911786
$ws.Range("D7").Value = 0.05245651392094138
$ws.Range("C8").Value = 0.8996581769135537
$ws.Range("D8").Value = 0.0224380259920453
$ws.Range("C9").Value = 0.8599938023002706
$ws.Range("D9").Value = 0.03240292840598669
$ws.Range("C10").Value = 0.8320840411457112
$ws.Range("D10").Value = 0.03568037894239338
$ws.Range("C11").Value = 0.7348810412912948
$ws.Range("D11").Value = 0.06081999668218552
$ws.Range("C12").Value = 0.9042680000000002
$ws.Range("D12").Value = 0.01957972128796234
$ws.Range("C13").Value = 0.8320840411457112
$ws.Range("D13").Value = 0.03568037894239338

$ws = $wb.Worksheets.Item("SVM")
$ws.Range("C2").Value = 0.7194075060401984
$ws.Range("D2").Value = 0.04687961694863366
$ws.Range("C3").Value = 0.8775187969924813
$ws.Range("D3").Value = 0.02052786939068589
$ws.Range("C4").Value = 0.9498015990799407
$ws.Range("D4").Value = 0.04139570896631687
$ws.Range("C5").Value = 0.5378878490790255
$ws.Range("D5").Value = 0.08115170617110701
$ws.Range("C6").Value = 0.990246
$ws.Range("D6").Value = 0.008868430825468641
$ws.Range("C7").Value = 0.6830611764451118
$ws.Range("D7").Value = 0.06889652147873214
$ws.Range("C8").Value = 0.8639790153571043
$ws.Range("D8").Value = 0.02585274267115796
$ws.Range("C9").Value = 0.8035261877487397
$ws.Range("D9").Value = 0.04017388151711887
$ws.Range("C10").Value = 0.7640663961264074
$ws.Range("D10").Value = 0.04038008586937003
$ws.Range("C11").Value = 0.6551108006471744
$ws.Range("D11").Value = 0.06349076942072276
$ws.Range("C12").Value = 0.8663099999999999
$ws.Range("D12").Value = 0.02012039527749018
$ws.Range("C13").Value = 0.7640663961264074
$ws.Range("D13").Value = 0.04038008586937002
